$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the cell values below can be updated.
$ws.Unprotect()

# Update the confidential disclaimer date from 2021-03-23 to 2021-03-24
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."
$ws.Rows(10).EntireRow.AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2491460688954086
$ws.Range("E2").Value = -0.0103213345039509

$ws.Range("D3").Value = 0.4942698923805518
$ws.Range("E3").Value = 0.002547410133031391

$ws.Range("D4").Value = 0.09946853439264516
$ws.Range("E4").Value = -0.01872201872201862

$ws.Range("D5").Value = 0.0997131600106079
$ws.Range("E5").Value = -0.004435245416912981

$ws.Range("D6").Value = 0.05740234432078643
$ws.Range("E6").Value = -0.01711433325410017

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = -0.004599318734477942
